$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 45484.184
$ws.Range("I11").Value = 45484.184
$ws.Range("K11").Value = 45484.184
$ws.Range("M11").Value = -45344.184

$ws.Range("H19").Value = 2143.889
$ws.Range("J19").Value = 2316.1667
$ws.Range("L19").Value = 2316.1667
$ws.Range("N19").Value = -2666.1667

$ws.Range("H32").Value = 7384.5
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H80").Value = 5488.125
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5488.125
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 16464.375
$ws.Range("N80").Value = -18460.375
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 5488.125
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5488.125
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 49393.125
$ws.Range("N83").Value = -59377.125
$ws.Range("M83").ClearContents()

$ws.Range("H92").Value = 306.9524
$ws.Range("I92").Value = 285.76923
$ws.Range("J92").Value = 341.375
$ws.Range("K92").Value = 285.76923
$ws.Range("L92").Value = 341.375
$ws.Range("M92").Value = 962.23077
$ws.Range("N92").Value = -2837.375

$ws.Range("H98").Value = 2089.5386
$ws.Range("I98").Value = 1742
$ws.Range("J98").Value = 3248
$ws.Range("K98").Value = 1742
$ws.Range("L98").Value = 3248
$ws.Range("M98").Value = -244
$ws.Range("N98").Value = -6244

$ws.Range("H122").Value = 2089.5386
$ws.Range("I122").Value = 1742
$ws.Range("J122").Value = 3248
$ws.Range("K122").Value = 5226
$ws.Range("L122").Value = 9744
$ws.Range("M122").Value = -2776
$ws.Range("N122").Value = -14644

$ws.Range("H132").Value = 58905.445
$ws.Range("I132").Value = 63956.188
$ws.Range("K132").Value = 191868.564
$ws.Range("M132").Value = -189338.564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3430.5454
$ws.Range("I2").Value = 2462.25
$ws.Range("J2").Value = 6012.6665
$ws.Range("K2").Value = 2462.25
$ws.Range("L2").Value = 6012.6665
$ws.Range("M2").Value = -2349.25
$ws.Range("N2").Value = -6238.6665

$ws.Range("H32").Value = 1022.3617
$ws.Range("I32").Value = 1033.8695
$ws.Range("J32").Value = 493
$ws.Range("K32").Value = 1033.8695
$ws.Range("L32").Value = 493
$ws.Range("M32").Value = -746.8695
$ws.Range("N32").Value = -1067

$ws.Range("H45").Value = 1122.4
$ws.Range("I45").Value = 1015.5
$ws.Range("J45").Value = 1550
$ws.Range("K45").Value = 1015.5
$ws.Range("L45").Value = 1550
$ws.Range("M45").Value = -638.5
$ws.Range("N45").Value = -2304

$ws.Range("H61").Value = 2385.353
$ws.Range("I61").Value = 1698.2727
$ws.Range("K61").Value = 1698.2727
$ws.Range("M61").Value = -1486.2727

$ws.Range("H74").Value = 2316448
$ws.Range("I74").Value = 1236295.5
$ws.Range("J74").Value = 5556905
$ws.Range("K74").Value = 1236295.5
$ws.Range("L74").Value = 5556905
$ws.Range("M74").Value = -1235421.5
$ws.Range("N74").Value = -5558653

$ws.Range("H77").Value = 2316448
$ws.Range("I77").Value = 1236295.5
$ws.Range("J77").Value = 5556905
$ws.Range("K77").Value = 6181477.5
$ws.Range("L77").Value = 27784525
$ws.Range("M77").Value = -6177109.5
$ws.Range("N77").Value = -27793261

$ws.Range("H102").Value = 926.7143
$ws.Range("I102").Value = 814.6667
$ws.Range("J102").Value = 1599
$ws.Range("K102").Value = 814.6667
$ws.Range("L102").Value = 1599
$ws.Range("M102").Value = 807.3333
$ws.Range("N102").Value = -4843

$ws.Range("H110").Value = 32364.4
$ws.Range("I110").Value = 29438.5
$ws.Range("J110").Value = 44068
$ws.Range("K110").Value = 29438.5
$ws.Range("L110").Value = 44068
$ws.Range("M110").Value = -27393.5
$ws.Range("N110").Value = -48158

$ws.Range("H116").Value = 3430.5454
$ws.Range("I116").Value = 2462.25
$ws.Range("J116").Value = 6012.6665
$ws.Range("K116").Value = 2462.25
$ws.Range("L116").Value = 6012.6665
$ws.Range("M116").Value = -168.25
$ws.Range("N116").Value = -10600.6665

$ws.Range("H122").Value = 1698.5238
$ws.Range("I122").Value = 1484.3077
$ws.Range("J122").Value = 2046.625
$ws.Range("K122").Value = 4452.9231
$ws.Range("L122").Value = 6139.875
$ws.Range("M122").Value = -2002.9231
$ws.Range("N122").Value = -11039.875

$ws.Range("H132").Value = 21741486
$ws.Range("I132").Value = 2378.6316
$ws.Range("J132").Value = 125002250
$ws.Range("K132").Value = 7135.8948
$ws.Range("L132").Value = 375006750
$ws.Range("M132").Value = -4605.8948
$ws.Range("N132").Value = -375011810

$ws.Range("H136").Value = 2385.353
$ws.Range("I136").Value = 1698.2727
$ws.Range("K136").Value = 5094.8181
$ws.Range("M136").Value = -2544.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3430.5454
$ws.Range("I3").Value = 2462.25
$ws.Range("J3").Value = 6012.6665
$ws.Range("K3").Value = 2462.25
$ws.Range("L3").Value = 6012.6665
$ws.Range("M3").Value = -2348.25
$ws.Range("N3").Value = -6240.6665

$ws.Range("H86").Value = 2495.1177
$ws.Range("I86").Value = 1833.4286
$ws.Range("J86").Value = 2958.3
$ws.Range("K86").Value = 1833.4286
$ws.Range("L86").Value = 2958.3
$ws.Range("M86").Value = -710.4286
$ws.Range("N86").Value = -5204.3

$ws.Range("H89").Value = 2495.1177
$ws.Range("I89").Value = 1833.4286
$ws.Range("J89").Value = 2958.3
$ws.Range("K89").Value = 9167.143
$ws.Range("L89").Value = 14791.5
$ws.Range("M89").Value = -3551.143
$ws.Range("N89").Value = -26023.5

$ws.Range("H94").Value = 1933.3684
$ws.Range("I94").Value = 1577.8823
$ws.Range("K94").Value = 1577.8823
$ws.Range("M94").Value = -1126.8823

$ws.Range("H99").Value = 2065.1428
$ws.Range("I99").Value = 2068.6956
$ws.Range("K99").Value = 2068.6956
$ws.Range("M99").Value = -570.6956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1789.591
$ws.Range("I31").Value = 1706.7894
$ws.Range("K31").Value = 1706.7894
$ws.Range("M31").Value = -1411.7894

$ws.Range("H34").Value = 1789.591
$ws.Range("I34").Value = 1706.7894
$ws.Range("K34").Value = 1706.7894
$ws.Range("M34").Value = -1504.7894

$ws.Range("H107").Value = 10252.23
$ws.Range("I107").Value = 14475.375
$ws.Range("J107").Value = 3495.2
$ws.Range("K107").Value = 14475.375
$ws.Range("L107").Value = 3495.2
$ws.Range("M107").Value = -12555.375
$ws.Range("N107").Value = -7335.2

$ws.Range("H132").Value = 5713.4546
$ws.Range("I132").Value = 5976.56
$ws.Range("J132").Value = 4891.25
$ws.Range("K132").Value = 17929.68
$ws.Range("L132").Value = 14673.75
$ws.Range("M132").Value = -15399.68
$ws.Range("N132").Value = -19733.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 12110.333
$ws.Range("J54").Value = 12110.333
$ws.Range("L54").Value = 36330.999
$ws.Range("N54").Value = -37448.999

$ws.Range("H116").Value = 94301.336
$ws.Range("I116").Value = 107347.695
$ws.Range("K116").Value = 322043.085
$ws.Range("M116").Value = -318601.085

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9894.611000000001
$ws.Range("I70").Value = 9078.929
$ws.Range("J70").Value = 12749.5
$ws.Range("K70").Value = 9078.929
$ws.Range("L70").Value = 12749.5
$ws.Range("M70").Value = -8808.929
$ws.Range("N70").Value = -13289.5

$ws.Range("H73").Value = 9894.611000000001
$ws.Range("I73").Value = 9078.929
$ws.Range("J73").Value = 12749.5
$ws.Range("K73").Value = 9078.929
$ws.Range("L73").Value = 12749.5
$ws.Range("M73").Value = -8142.929
$ws.Range("N73").Value = -14621.5

$ws.Range("H80").Value = 4919.0454
$ws.Range("I80").Value = 3155.3125
$ws.Range("J80").Value = 9622.333000000001
$ws.Range("K80").Value = 3155.3125
$ws.Range("L80").Value = 9622.333000000001
$ws.Range("M80").Value = -2157.3125
$ws.Range("N80").Value = -11618.333

$ws.Range("H83").Value = 4919.0454
$ws.Range("I83").Value = 3155.3125
$ws.Range("J83").Value = 9622.333000000001
$ws.Range("K83").Value = 15776.5625
$ws.Range("L83").Value = 48111.665
$ws.Range("M83").Value = -10784.5625
$ws.Range("N83").Value = -58095.665

$ws.Range("H107").Value = 691.4
$ws.Range("I107").Value = 528.58826
$ws.Range("K107").Value = 528.58826
$ws.Range("M107").Value = 1391.41174

$ws.Range("H132").Value = 909.5
$ws.Range("I132").Value = 855.4
$ws.Range("K132").Value = 2566.2
$ws.Range("M132").Value = -36.19999999999982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1857.8
$ws.Range("I7").Value = 1556.6
$ws.Range("K7").Value = 1556.6
$ws.Range("M7").Value = -1444.6

$ws.Range("H68").Value = 2926.6667
$ws.Range("I68").Value = 2926.6667
$ws.Range("K68").Value = 2926.6667
$ws.Range("M68").Value = -2177.6667

$ws.Range("H71").Value = 2926.6667
$ws.Range("I71").Value = 2926.6667
$ws.Range("K71").Value = 14633.3335
$ws.Range("M71").Value = -10889.3335

$ws.Range("H126").Value = 1857.8
$ws.Range("I126").Value = 1556.6
$ws.Range("K126").Value = 4669.799999999999
$ws.Range("M126").Value = -2199.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3130.5833
$ws.Range("I96").Value = 4515.5454
$ws.Range("J96").Value = 1958.6923
$ws.Range("K96").Value = 4515.5454
$ws.Range("L96").Value = 1958.6923
$ws.Range("M96").Value = -3142.5454
$ws.Range("N96").Value = -4704.6923
